# The commit (an "Add files via upload" re-save of the workbook from a
# different machine/build of Excel) renames the single worksheet from
# "Planilha1" to "Dados1". That rename is the only deliberate, content-level
# edit captured in the diff -- everything else (fileVersion/rupBuild,
# x15ac:absPath, xr:revisionPtr GUID, workbookView window geometry,
# x14ac:dyDescent hints, the tiny column-width deltas, and the customXml
# SharePoint content-type part renumbering) is environment/build metadata
# that Excel itself regenerates on every save and isn't something the
# workbook's object model exposes for a script to set.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Name = "Dados1"
